$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Data rows 3-30 (rows 1-2 are unchanged header rows).
# A new row "NIY" is inserted at row 3, shifting the previously
# existing rows down by one; several new rows are appended at the
# end (14, 16, 20, 23-30 new content) implementing the extended
# sign-up / military-invest / government event flow.
# ------------------------------------------------------------------

$ws.Range("A3").Value = 'NIY'
$ws.Range("B3").Value = '功能还没有实现'
$ws.Range("C3").Value = 'dialog'
$ws.Range("D3").Value = 'dialog_not_implement_yet'

$ws.Range("A4").Value = 'government'
$ws.Range("B4").Value = '皇宫'
$ws.Range("C4").Value = 'selectlist'
$ws.Range("D4").Value = 'signUp;militaryInvest;recommendLetter;close'

$ws.Range("A5").Value = 'tarven'
$ws.Range("B5").Value = '酒馆'
$ws.Range("C5").Value = 'selectlist'
$ws.Range("D5").Value = 'drink;buyDrink;hireSailor;spreadRumor;close'

$ws.Range("A6").Value = 'plaza'
$ws.Range("B6").Value = '广场'
$ws.Range("C6").Value = 'plazaEvent'
$ws.Range("D6").Value = ';'

$ws.Range("A7").Value = 'exchange'
$ws.Range("B7").Value = '交易所'
$ws.Range("C7").Value = 'selectlist'
$ws.Range("D7").Value = 'trade;commerceInvest;commerceInfo;close'

$ws.Range("A8").Value = 'shipyard'
$ws.Range("B8").Value = '造船厂'
$ws.Range("C8").Value = 'selectlist'
$ws.Range("D8").Value = 'buyShip;sellShip;modifyShip;fixShip;dockYard;close'

$ws.Range("A9").Value = 'shop'
$ws.Range("B9").Value = '商店'
$ws.Range("C9").Value = 'selectlist'
$ws.Range("D9").Value = 'buyItem;sellItem;taskList;mailLetter;close'

$ws.Range("A10").Value = 'dock'
$ws.Range("B10").Value = '船坞'
$ws.Range("C10").Value = 'selectlist'
$ws.Range("D10").Value = 'sail;supply;dockYard;close'

$ws.Range("A11").Value = 'inn'
$ws.Range("B11").Value = '旅馆'
$ws.Range("C11").Value = 'selectlist'
$ws.Range("D11").Value = 'sleep1day;sleep10day;sleep30day;close'

$ws.Range("A12").Value = 'relic'
$ws.Range("B12").Value = '遗迹'
$ws.Range("C12").Value = 'relicEvent'
$ws.Range("D12").Value = ';'

$ws.Range("A13").Value = 'close'
$ws.Range("B13").Value = '关闭窗口'
$ws.Range("C13").Value = 'close'
$ws.Range("D13").Value = ';'

$ws.Range("A14").Value = 'waitADay'
$ws.Range("B14").Value = '等待一天'
$ws.Range("C14").Value = 'wait'
$ws.Range("D14").Value = 1

$ws.Range("A15").Value = 'canSignUp'
$ws.Range("B15").Value = '条件分歧'
$ws.Range("C15").Value = 'condition'
$ws.Range("D15").Value = 'canSignUpMoneyEnough;signUpFailed'

$ws.Range("A16").Value = 'canSignUpMoneyEnough'
$ws.Range("B16").Value = '条件分歧'
$ws.Range("C16").Value = 'condition'
$ws.Range("D16").Value = 'signUpStart;signUpMoneyNotEnough'

$ws.Range("A17").Value = 'signUpMoneyNotEnough'
$ws.Range("B17").Value = '签约钱不够'
$ws.Range("C17").Value = 'dialog'
$ws.Range("D17").Value = 'dialog_no_enough_money'

$ws.Range("A18").Value = 'signUpStart'
$ws.Range("B18").Value = '签约开始'
$ws.Range("C18").Value = 'eventList'
$ws.Range("D18").Value = 'signUpDialog;signUpWindow;signUpClose'

$ws.Range("A19").Value = 'signUpDialog'
$ws.Range("B19").Value = '签约对话'
$ws.Range("C19").Value = 'dialog'
$ws.Range("D19").Value = 'dialog_signup_start'

$ws.Range("A20").Value = 'signUpWindow'
$ws.Range("B20").Value = '签约窗口'
$ws.Range("C20").Value = 'window'
$ws.Range("D20").Value = 'InvestPanel;2;signUpSuccess;signUpFail'

$ws.Range("A21").Value = 'signUpFailed'
$ws.Range("B21").Value = '签约失败对话'
$ws.Range("C21").Value = 'dialog'
$ws.Range("D21").Value = 'dialog_signup_failure_full'

$ws.Range("A22").Value = 'signUp'
$ws.Range("B22").Value = '签约'
$ws.Range("C22").Value = 'eventList'
$ws.Range("D22").Value = 'canSignUp'

$ws.Range("A23").Value = 'signUpSuccess'
$ws.Range("B23").Value = '签约成功'
$ws.Range("C23").Value = 'eventList'
$ws.Range("D23").Value = 'signUpSuccessDialog;close;waitADay;government'

$ws.Range("A24").Value = 'signUpSuccessDialog'
$ws.Range("B24").Value = '签约成功对话'
$ws.Range("C24").Value = 'dialog'
$ws.Range("D24").Value = 'dialog_invest_success'

$ws.Range("A25").Value = 'signUpFail'
$ws.Range("B25").Value = '签约失败对话'
$ws.Range("C25").Value = 'dialog'
$ws.Range("D25").Value = 'dialog_invest_fail'

$ws.Range("A26").Value = 'wannaInvestMilitary'
$ws.Range("B26").Value = '想投资军事'
$ws.Range("C26").Value = 'dialog'
$ws.Range("D26").Value = 'dalog_wanna_invest'

$ws.Range("A27").Value = 'militaryInvest'
$ws.Range("B27").Value = '军事投资'
$ws.Range("C27").Value = 'eventList'
$ws.Range("D27").Value = 'wannaInvestMilitary;canMilitaryInvestMoneyEnough'

$ws.Range("A28").Value = 'canMilitaryInvestMoneyEnough'
$ws.Range("B28").Value = '条件分歧'
$ws.Range("C28").Value = 'condition'
$ws.Range("D28").Value = 'militaryInvestStart;signUpMoneyNotEnough'

$ws.Range("A29").Value = 'militaryInvestStart'
$ws.Range("B29").Value = '军事投资开始'
$ws.Range("C29").Value = 'eventList'
$ws.Range("D29").Value = 'militaryWindow;signUpClose'

$ws.Range("A30").Value = 'militaryWindow'
$ws.Range("B30").Value = '军事投资窗口'
$ws.Range("C30").Value = 'window'
$ws.Range("D30").Value = 'InvestPanel;1;signUpSuccess;signUpFail'

# ------------------------------------------------------------------
# Re-apply the "CJK body text" cell style (font: SimSun / 宋体) to the
# description/eventType columns, matching the look of the original
# rows. Copy-format from an already-styled source cell so the xlsx
# reuses the existing style record instead of minting a new one.
# ------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B4:B30").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C4:C14").PasteSpecial(-4122)
$ws.Range("C17:C27").PasteSpecial(-4122)
$ws.Range("C29:C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Column widths (A gains a bestFit width, a new width is set for the
# newly-visible column C, B/D widen slightly).
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.330729166666666
$ws.Columns.Item(3).ColumnWidth = 11.330729166666666
$ws.Columns.Item(4).ColumnWidth = 43.330729166666664

# ------------------------------------------------------------------
# Selection / view: author scrolled down and left the cursor on D29.
# ------------------------------------------------------------------
$ws.Range("D29").Select()

Write-Host "eventAction sheet updated"